$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 496 (shifts existing rows 496:559 down to 497:560)
$ws.Rows.Item(496).Insert()

# Populate the newly inserted row with the new "Mango" price-record data
$ws.Range("A496").Value = 9
$ws.Range("B496").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C496").Value = "Metropolitana"
$ws.Range("D496").Value = 44918
$ws.Range("E496").Value = 13
$ws.Range("F496").Value = "Fruta"
$ws.Range("G496").Value = 100108
$ws.Range("H496").Value = "Tropicales y subtropicales"
$ws.Range("I496").Value = 100108002
$ws.Range("J496").Value = "Mango"
$ws.Range("K496").Value = "Sin especificar"
$ws.Range("L496").Value = "Primera"
$ws.Range("M496").Value = 680
$ws.Range("N496").Value = 6500
$ws.Range("O496").Value = 7500
$ws.Range("P496").Value = 7059
$ws.Range("Q496").Value = '$/bandeja 4 kilos'
$ws.Range("R496").Value = "Perú"
$ws.Range("S496").Value = 1765
$ws.Range("T496").Value = 4
